$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Juan) typed across the whole row first.
$ws.Range("A3").Value = "Juan"
$ws.Range("B3").Value = "juan@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:juan@gmail.com")
$ws.Range("B3").Style = "Hipervínculo"
$ws.Range("C3").Value = "75643234W"
$ws.Range("D3").Value = "jkcfdasd"

# Remaining six rows (4-9) filled column by column: all names, then all
# emails, then all nifs, then all codes - matching how the shared-string
# table grew in the source workbook.
$names   = @("Sergio", "David", "Santiago", "Xurso", "Adrian", "Luis")
$emails  = @("sergio@yomolomucho.es", "david@gmail.com", "santi@yomolomucho.es", "xurso@gmail.com", "adrian@yomolomucho.es", "luis@gmail.com")
$nifs    = @("12321543P", "89098456D", "12047623S", "71234432X", "56412376R", "54234981Q")
$codigos = @("fdshghg", "jkcfdasd", "fdshghg", "jkcfdasd", "fdshghg", "jkcfdasd")

for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 4
    $ws.Range("A$r").Value = $names[$i]
}
for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 4
    $ws.Range("B$r").Value = $emails[$i]
    $ws.Hyperlinks.Add($ws.Range("B$r"), "mailto:" + $emails[$i])
    $ws.Range("B$r").Style = "Hipervínculo"
}
for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 4
    $ws.Range("C$r").Value = $nifs[$i]
}
for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 4
    $ws.Range("D$r").Value = $codigos[$i]
}

# Page setup tweak picked up by the saved workbook.
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

# Leave the selection where the user ended up after entering the data.
$ws.Range("C10").Select()
